# Apply the HEVDoc.xlsx content update:
# - Refresh the "Chapters" sheet text (trademark symbols, reworded abstract, reordered entries)
# - Update window/selection state

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chapters")

# Write in the same order the new shared strings were introduced upstream
# (report summary, toolbox list w/ trademarks, release w/ trademark, reworded abstract)
$ws.Range("B2").Value = "The report summarizes the results of a fleet test campaign of the hybrid electric vehicle (HEV) P4 model. The document provides both the model setup info and results.`n"
$ws.Range("B4").Value = "MATLAB®`nSimulink®`nPowertrain Blockset™`nMATLAB® Report Generator™`nSimulink® Compiler™`nMATLAB® Compiler™`nMATLAB® Web App Server™"
$ws.Range("B3").Value = "MATLAB® R2021a"
$ws.Range("B5").Value = "The Hybrid Electric Vehicle (HEV) P4 Reference Application represents a full HEV model with an internal combustion engine, transmission, battery, motor, and associated powertrain control algorithms. Use the reference application for hardware-in-the-loop (HIL) testing, tradeoff analysis, and control parameter optimization of a HEV P4 hybrid. "

$ws.Activate() | Out-Null
$ws.Range("B6").Select() | Out-Null
